$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data stores every Coin/Link/Price/Volume cell as literal
# text (Coinranking renders "88.003.57" style thousand-separated prices,
# not real numbers). Any replacement value that Excel would otherwise
# auto-detect as a number (e.g. "0.130", "13.50") is written to a cell
# pre-formatted as Text so it is stored verbatim, matching the original
# authoring (equivalent to typing it with a leading apostrophe).

$ws.Range("D2").Value = "87.610.34"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "3.231.96"
$ws.Range("E3").Value = "  -2.94%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.09"
$ws.Range("E5").Value = "  -5.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "613.04"
$ws.Range("E6").Value = "  -5.98%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.384"
$ws.Range("E7").Value = "  +7.90%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.680"
$ws.Range("E8").Value = "  +12.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").Value = "3.225.84"
$ws.Range("E10").Value = "  -3.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.545"
$ws.Range("E11").Value = "  -7.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.181"
$ws.Range("E12").Value = "  +7.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000249"
$ws.Range("E13").Value = "  -7.88%  "
$ws.Range("D14").Value = "3.819.83"
$ws.Range("E14").Value = "  -3.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.33"
$ws.Range("E15").Value = "  -3.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "32.77"
$ws.Range("E16").Value = "  -7.55%  "
$ws.Range("D17").Value = "87.483.11"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("D18").Value = "3.232.27"
$ws.Range("E18").Value = "  -2.70%  "
$ws.Range("B19").Value = "SuiNetwork"
$ws.Range("C19").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.98"
$ws.Range("E19").Value = "  -5.45%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.50"
$ws.Range("E20").Value = "  -8.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "422.33"
$ws.Range("E21").Value = "  -7.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.61"
$ws.Range("E22").Value = "  -12.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.14"
$ws.Range("E23").Value = "  -7.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.22"
$ws.Range("E24").Value = "  -6.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.69"
$ws.Range("E25").Value = "  -8.17%  "
$ws.Range("D26").Value = "3.388.16"
$ws.Range("E26").Value = "  -2.87%  "
$ws.Range("B27").Value = "Litecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "74.66"
$ws.Range("E27").Value = "  -5.60%  "
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000133"
$ws.Range("E28").Value = "  +5.33%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.174"
$ws.Range("E30").Value = "  -11.41%  "
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "549.19"
$ws.Range("E32").Value = "  -9.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.48"
$ws.Range("E33").Value = "  -10.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.89"
$ws.Range("E34").Value = "  -10.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.28"
$ws.Range("E35").Value = "  -20.54%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.73"
$ws.Range("E36").Value = "  -6.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.135"
$ws.Range("E37").Value = "  -7.95%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "22.38"
$ws.Range("E38").Value = "  -4.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "21.84"
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.01"
$ws.Range("E41").Value = "  -0.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.385"
$ws.Range("E42").Value = "  -8.74%  "
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.91"
$ws.Range("E44").Value = "  -11.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "146.88"
$ws.Range("E45").Value = "  -7.95%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "175.02"
$ws.Range("E46").Value = "  -8.79%  "
$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "43.78"
$ws.Range("E47").Value = "  -5.27%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.130"
$ws.Range("E48").Value = "  +13.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.29"
$ws.Range("E49").Value = "  -9.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.07"
$ws.Range("E50").Value = "  -9.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.607"
$ws.Range("E51").Value = "  -8.54%  "
